$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Character used for the "horizontal ellipsis" glyph in the WordPress "(more...)" teaser text.
$ellipsis = [char]0x2026

# --- Row 3: Shark Tank -------------------------------------------------
$ws.Range("A3").Value = "Shark Tank"
$ws.Range("B3").Value = "Episode 904"
$ws.Range("C3").Value = "(Season 9, Episode 2)"
$ws.Range("D3").Value = "Rohan Oza sits in with the sharks. Products include ultimate survival kits; a wireless microphone designed to engage kids; and an app to navigate airports."
$ws.Range("E3").Value = "ABC"
# Date column keeps being stored as text (not an auto-converted date serial),
# same as the existing Text number format already used on this column.
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "Dec 3"
$ws.Range("G3").Value = "10:00pm"

# --- Row 4: Star Trek: Voyager ------------------------------------------
$ws.Range("A4").Value = "Star Trek: Voyager"
$ws.Range("B4").Value = "Unimatrix Zero"
$ws.Range("C4").Value = "(Season 6, Episode 26)"
$ws.Range("D4").Value = "Part 1 of two. Seven of Nine is drawn into a dreamworld that Borg drones inhabit during their sleep cycles---a threat to the Borg Queen's control that Janeway wants to exploit. Borg Queen: Susanna Thompson. Axum: Mark Deakins. Korok: Jerome Butler. Laura: (more$ellipsis)Part 1 of two. Seven of Nine is drawn into a dreamworld that Borg drones inhabit during their sleep cycles---a threat to the Borg Queen's control that Janeway wants to exploit. Borg Queen: Susanna Thompson. Axum: Mark Deakins. Korok: Jerome Butler. Laura: Joanna Heimbold. Seven of Nine: Jeri Ryan. Janeway: Kate Mulgrew."
$ws.Range("E4").Value = "BBC"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "Dec 4"
$ws.Range("G4").Value = "8:00pm"
